$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (Turkish 2 Lig); remaining rows shift up
$ws.Rows("8").Delete() | Out-Null

# Pre-format Date/Time columns (B and C) as Text so values are not converted to date/time serials
$dtRange = $ws.Range("B2:C7")
$dtRange.NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'Portuguese Segunda Liga'
$ws.Range("B2").Value = '2025-12-16'
$ws.Range("C2").Value = '14:00:00'
$ws.Range("D2").Value = 'Maritimo'
$ws.Range("E2").Value = 'Benfica B'
$ws.Range("F2").Value = 1.79
$ws.Range("G2").Value = 1.81
$ws.Range("H2").Value = 4.9
$ws.Range("I2").Value = 5.1
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.2
$ws.Range("L2").Value = 1.52
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.85
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 1.91
$ws.Range("Q2").Value = 2.06
$ws.Range("R2").Value = 1.34
$ws.Range("S2").Value = 3.8
$ws.Range("T2").Value = 1.05
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 16
$ws.Range("Y2").Value = 970
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 970
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 480
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 480
$ws.Range("AH2").Value = 490
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 24
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("A3").Value = 'Swiss Super League'
$ws.Range("B3").Value = '2025-12-16'
$ws.Range("C3").Value = '16:30:00'
$ws.Range("D3").Value = 'St Gallen'
$ws.Range("E3").Value = 'Sion'
$ws.Range("F3").Value = 2.58
$ws.Range("G3").Value = 2.64
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 2.94
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 3.7
$ws.Range("L3").Value = 1.37
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 4.2
$ws.Range("O3").Value = 1.29
$ws.Range("P3").Value = 2.1
$ws.Range("Q3").Value = 1.89
$ws.Range("R3").Value = 1.43
$ws.Range("S3").Value = 3.2
$ws.Range("T3").Value = 1.69
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.51
$ws.Range("W3").Value = 1.61
$ws.Range("X3").Value = 17.5
$ws.Range("Y3").Value = 13
$ws.Range("Z3").Value = 20
$ws.Range("AA3").Value = 46
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 8
$ws.Range("AD3").Value = 12.5
$ws.Range("AE3").Value = 32
$ws.Range("AF3").Value = 17.5
$ws.Range("AG3").Value = 11.5
$ws.Range("AH3").Value = 16.5
$ws.Range("AI3").Value = 40
$ws.Range("AJ3").Value = 38
$ws.Range("AK3").Value = 27
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 80
$ws.Range("AN3").Value = 19.5
$ws.Range("AO3").Value = 25

# Row 4
$ws.Range("A4").Value = 'Swiss Super League'
$ws.Range("B4").Value = '2025-12-16'
$ws.Range("C4").Value = '16:30:00'
$ws.Range("D4").Value = 'Winterthur'
$ws.Range("E4").Value = 'Thun'
$ws.Range("F4").Value = 4.5
$ws.Range("G4").Value = 4.7
$ws.Range("H4").Value = 1.75
$ws.Range("I4").Value = 1.76
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 4.7
$ws.Range("L4").Value = 1.28
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 6.4
$ws.Range("O4").Value = 1.17
$ws.Range("P4").Value = 2.84
$ws.Range("Q4").Value = 1.52
$ws.Range("R4").Value = 1.75
$ws.Range("S4").Value = 2.28
$ws.Range("T4").Value = 1.57
$ws.Range("U4").Value = 2.6
$ws.Range("V4").Value = 2.3
$ws.Range("W4").Value = 1.27
$ws.Range("X4").Value = 29
$ws.Range("Y4").Value = 14
$ws.Range("Z4").Value = 14
$ws.Range("AA4").Value = 19.5
$ws.Range("AB4").Value = 27
$ws.Range("AC4").Value = 11.5
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 15.5
$ws.Range("AF4").Value = 40
$ws.Range("AG4").Value = 18
$ws.Range("AH4").Value = 16.5
$ws.Range("AI4").Value = 25
$ws.Range("AJ4").Value = 100
$ws.Range("AK4").Value = 46
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 60
$ws.Range("AN4").Value = 34
$ws.Range("AO4").Value = 6.8

# Row 5
$ws.Range("A5").Value = 'English National League'
$ws.Range("B5").Value = '2025-12-16'
$ws.Range("C5").Value = '16:45:00'
$ws.Range("D5").Value = 'Truro City'
$ws.Range("E5").Value = 'Wealdstone'
$ws.Range("F5").Value = 3.1
$ws.Range("G5").Value = 3.2
$ws.Range("H5").Value = 2.42
$ws.Range("I5").Value = 2.54
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 3.65
$ws.Range("L5").Value = 1.43
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.75
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 3.6
$ws.Range("T5").Value = 1.75
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.65
$ws.Range("W5").Value = 1.45
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 16
$ws.Range("AA5").Value = 36
$ws.Range("AB5").Value = 12.5
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 27
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 13.5
$ws.Range("AH5").Value = 17
$ws.Range("AI5").Value = 40
$ws.Range("AJ5").Value = 60
$ws.Range("AK5").Value = 36
$ws.Range("AL5").Value = 46
$ws.Range("AM5").Value = 95
$ws.Range("AN5").Value = 42
$ws.Range("AO5").Value = 23

# Row 6
$ws.Range("A6").Value = 'Welsh Premiership'
$ws.Range("B6").Value = '2025-12-16'
$ws.Range("C6").Value = '16:45:00'
$ws.Range("D6").Value = 'Cardiff Metropolitan'
$ws.Range("E6").Value = 'Briton Ferry Llansawel'
$ws.Range("F6").Value = 1.98
$ws.Range("G6").Value = 2.06
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 3.95
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 4.3
$ws.Range("L6").Value = 1.35
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 5.1
$ws.Range("O6").Value = 1.22
$ws.Range("P6").Value = 2.34
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 1.54
$ws.Range("S6").Value = 2.74
$ws.Range("T6").Value = 1.66
$ws.Range("U6").Value = 2.4
$ws.Range("V6").Value = 1.33
$ws.Range("W6").Value = 1.94
$ws.Range("X6").Value = 19.5
$ws.Range("Y6").Value = 19
$ws.Range("Z6").Value = 32
$ws.Range("AA6").Value = 80
$ws.Range("AB6").Value = 13
$ws.Range("AC6").Value = 9.4
$ws.Range("AD6").Value = 16
$ws.Range("AE6").Value = 130
$ws.Range("AF6").Value = 15.5
$ws.Range("AG6").Value = 11
$ws.Range("AH6").Value = 16.5
$ws.Range("AI6").Value = 44
$ws.Range("AJ6").Value = 24
$ws.Range("AK6").Value = 19
$ws.Range("AL6").Value = 30
$ws.Range("AM6").Value = 70
$ws.Range("AN6").Value = 11
$ws.Range("AO6").Value = 32

# Row 7
$ws.Range("A7").Value = 'Colombian Primera A'
$ws.Range("B7").Value = '2025-12-16'
$ws.Range("C7").Value = '21:30:00'
$ws.Range("D7").Value = 'Tolima'
$ws.Range("E7").Value = 'Junior FC Barranquilla'
$ws.Range("F7").Value = 1.76
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 5.3
$ws.Range("I7").Value = 5.8
$ws.Range("J7").Value = 3.7
$ws.Range("K7").Value = 3.9
$ws.Range("L7").Value = 1.45
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 3.3
$ws.Range("O7").Value = 1.39
$ws.Range("P7").Value = 1.81
$ws.Range("Q7").Value = 2.18
$ws.Range("R7").Value = 1.29
$ws.Range("S7").Value = 4.1
$ws.Range("T7").Value = 2
$ws.Range("U7").Value = 1.85
$ws.Range("V7").Value = 1.21
$ws.Range("W7").Value = 2.24
$ws.Range("X7").Value = 13.5
$ws.Range("Y7").Value = 16.5
$ws.Range("Z7").Value = 42
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 7.6
$ws.Range("AC7").Value = 8.6
$ws.Range("AD7").Value = 22
$ws.Range("AE7").Value = 90
$ws.Range("AF7").Value = 10
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 24
$ws.Range("AI7").Value = 100
$ws.Range("AJ7").Value = 18.5
$ws.Range("AK7").Value = 20
$ws.Range("AL7").Value = 48
$ws.Range("AM7").Value = 160
$ws.Range("AN7").Value = 15
$ws.Range("AO7").Value = 130

# Restore default style for the text-formatted Date/Time cells (keep values as text, drop number formatting)
$dtRange.Style = "Normal"
